# zootable/test_data/example.xlsx
#
# "convert accession nums to char / ignore case on default ordering /
#  validate that all accession numbers are 6 digits"
#
# The accession numbers in column B (rows 2-5) were mistakenly entered with
# an extra leading digit (7 digits instead of the required 6). Trim them
# down to 6 digits each so they match the validation rule (and the existing
# 6-digit value already present in B6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 111111
$ws.Range("B3").Value = 111112
$ws.Range("B4").Value = 111113
$ws.Range("B5").Value = 111114

# Restore the sheet's default column width (harmless if the host ignores it).
$ws.StandardWidth = 11.53515625

# Move/record the active selection on the sheet, as it was when the author
# saved the file after making the edit above.
$ws.Range("B3").Select()
